# "ultimas cositas de las diapositivas :p"
#
# 1) Slide 5 ("Event C#"): nudge the title placeholder down a touch
#    (y offset 348865 -> 378362 EMU, i.e. 27.469685pt -> 29.792283pt).
# 2) Slide 7: title becomes "3. ¿Qué es un singleton?" (numbering the
#    question + splitting "singleton" into its own run).
# 3) Slide 9: title gains a leading space -> " Implementación en C# y Unity".

$p = $ppt.ActivePresentation

# --- Slide 5: reposition title -----------------------------------------
$s5 = $p.Slides.Item(5)
$title5 = $s5.Shapes.Title
$title5.Top = 378362 / 12700

# --- Slide 7: "¿Qué es un singleton?" -> "3. ¿Qué es un singleton?" ----
$s7 = $p.Slides.Item(7)
$title7 = $s7.Shapes.Title
$tr7 = $title7.TextFrame.TextRange
$tr7.Text = "3. ¿Qué es un singleton?"
# Re-split "singleton" and the trailing "?" into their own runs.
$tr7.Characters(15, 9).Text = "singleton"
$tr7.Characters(24, 1).Text = "?"

# --- Slide 9: "Implementación en C# y Unity" -> " Implementación en C# y Unity"
$s9 = $p.Slides.Item(9)
$title9 = $s9.Shapes.Title
$title9.TextFrame.TextRange.Text = " Implementación en C# y Unity"
